$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing metric values (rows 2-8) for the new Deep Learning comparison data ---
$ws.Range("C2").Value = 0.49
$ws.Range("E2").Value = 0.626
$ws.Range("F2").Value = 0.479

$ws.Range("C3").Value = 0.68
$ws.Range("E3").Value = 0.711
$ws.Range("F3").Value = 0.75

$ws.Range("C4").Value = 0.64
$ws.Range("D4").Value = 0.892
$ws.Range("E4").Value = 0.646
$ws.Range("F4").Value = 0.75

$ws.Range("C5").Value = 0.655
$ws.Range("E5").Value = 0.728
$ws.Range("F5").Value = 0.706

$ws.Range("C6").Value = 0.685
$ws.Range("E6").Value = 0.681
$ws.Range("F6").Value = 0.775

$ws.Range("C7").Value = 0.63
$ws.Range("E7").Value = 0.711
$ws.Range("F7").Value = 0.681

$ws.Range("C8").Value = 0.66
$ws.Range("E8").Value = 0.67
$ws.Range("F8").Value = 0.753

# --- Widen column B to fit the new "Deep Learning Models" header ---
$ws.Columns.Item(2).ColumnWidth = 27.5

# --- New section header row (row 10): "Deep Learning Models" styled like the title row ---
$ws.Range("B10").Value = "Deep Learning Models"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(10).RowHeight = 18.75

# --- New data row (row 11): ANN results ---
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "ANN"
$ws.Range("C11").Value = 0.605
$ws.Range("D11").Value = 0.644
$ws.Range("E11").Value = 0.684
$ws.Range("F11").Value = 0.663

# --- Clear the clipboard marching ants / update selection to match the final cursor position ---
$excel.CutCopyMode = $false
$ws.Range("F10").Select() | Out-Null
